# Apply conservative present-scenario LCOH values (capacity-factor adjusted)
# across the Present-Storage, 2030-Storage and 2050-Storage sheets.

$wb = $excel.ActiveWorkbook

# --- Present-Storage sheet: all LCOH values replaced, column B narrowed ---
$wsPresent = $wb.Worksheets.Item("Present-Storage")

$wsPresent.Range("B2").Value = 7.71
$wsPresent.Range("B3").Value = 5.53
$wsPresent.Range("B4").Value = 8.33
$wsPresent.Range("B5").Value = 8.24
$wsPresent.Range("B6").Value = 7
$wsPresent.Range("B7").Value = 4.7
$wsPresent.Range("B8").Value = 7.66
$wsPresent.Range("B9").Value = 7.56
$wsPresent.Range("B10").Value = 15.8
$wsPresent.Range("B11").Value = 14.11
$wsPresent.Range("B12").Value = 16.29
$wsPresent.Range("B13").Value = 16.22
$wsPresent.Range("B14").Value = 7.95
$wsPresent.Range("B15").Value = 5.86
$wsPresent.Range("B16").Value = 8.539999999999999
$wsPresent.Range("B17").Value = 8.449999999999999

$wsPresent.Columns.Item(2).ColumnWidth = 8.4

# --- 2030-Storage sheet: WindOnshore rows updated for capacity factor ---
$ws2030 = $wb.Worksheets.Item("2030-Storage")

$ws2030.Range("B3").Value = 4.41
$ws2030.Range("B7").Value = 4.14
$ws2030.Range("B11").Value = 8.41
$ws2030.Range("B15").Value = 5.55

# --- 2050-Storage sheet: WindOnshore rows updated for capacity factor ---
$ws2050 = $wb.Worksheets.Item("2050-Storage")

$ws2050.Range("B3").Value = 68.48999999999999
$ws2050.Range("B7").Value = 71.90000000000001
$ws2050.Range("B11").Value = 58.41
$ws2050.Range("B15").Value = 67.03
